# "updated checklist, added a few new items"
# - Tasks!D8: mark the "Writeup Hybrid algorithm" status as DONE (was 50%).
# - Tasks!B24:B26: three new checklist items appended below "What is the
#   .bgs file used for? ..." (row 23).
# - Move the selection down to the newly added rows, roughly matching the
#   author's final cursor position/view.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tasks")

$ws.Range("D8").Value = "DONE"

$ws.Range("B24").Value = "Add section about seeds/reproducibility"
$ws.Range("B25").Value = "Mention capabilities of R2admb"
$ws.Range("B26").Value = "define what mcmult argument does"

# Scroll/select close to the new content, like the author would have left
# the sheet after typing the new rows.
$ws.Range("D27").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 7
